$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C, rows 2 through 501: change date serial value from 45182 to 45184
for ($r = 2; $r -le 501; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
